# Applies the autofilter + row-visibility + two capacity-value updates
# described by the commit "Added H2 demand NOS0, updated VRE, updated demand NO."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Capacity")

# --- Data updates -----------------------------------------------------
# NOS0 Electrolysis (row 64) and NOS0 Hydrogen processor (row 66):
# Other_capa 500 -> 1500
$ws.Range("H64").Value = 1500
$ws.Range("H66").Value = 1500

# --- AutoFilter: restrict the Capacity table to Year 2040 and the
# Norwegian nodes NOM1 / NON1 / NOS0. Excel recomputes row visibility
# for the filtered range automatically.
$rng = $ws.Range("A1:J91")
$rng.AutoFilter(5, @("2040"))
$rng.AutoFilter(1, @("NOM1", "NON1", "NOS0"))

# --- Selection / scroll position shown in the diff ---------------------
$ws.Range("E98").Select()
